$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ducks")

# --- Duck name change: row 29 (Name column B), "The Donald" -> "Backquaker" ---
$ws.Range("B29").Value = 'Backquaker'

# --- Row 10: rewrite the Halloween group fact (now plural "our/we/we'll") and drop bold ---
$ws.Range("L10").Value = 'Trick or treat, smell our beaks, give us something good to eat, if you don''t we won''t care, we''ll just pull up your underwear.'
$ws.Range("L10").Font.Bold = $false

# --- New "About Me/Fun Fact" text for rows that previously held the placeholder --- 
$ws.Range("L14").Value = 'My hat''s propeller makes my migrations a kaleidoscope of color!'
$ws.Rows.Item(14).RowHeight = 28.8

$ws.Range("L18").Value = 'I''m 50% mermaid, 50% duck, and 100% that billtch'
$ws.Rows.Item(18).RowHeight = 28.8

$ws.Range("L26").Value = 'I want to become a Pokemon like my big bro Porygon. My signature move is Printed Peck Attack.'
$ws.Rows.Item(26).RowHeight = 43.2

$ws.Range("L27").Value = 'Don''t tell Jack, but I''m the real pumpkin king.'

$ws.Range("L30").Value = "If I were a rich duck.`nYa ba dibba dibba dibba dibba dibba dibba dum`nAll day long, I'd biddy biddy bum.`nIf I were a wealthy duck`nI wouldn't have to fly hard`nYa ba dibba dibba dibba dibba dibba dibba dum.`nIf I were a biddy biddy rich yidle-diddle-didle-didle duck"
$ws.Rows.Item(30).RowHeight = 129.6

$ws.Range("L31").Value = 'Finding the Ark of the Covenant and the Holy Quail is all in a day''s work for me.'
$ws.Rows.Item(31).RowHeight = 28.8

$ws.Range("L40").Value = 'I love long walks in the Fowlbidden Forest, one time I even saw Ron Webbsley!'
$ws.Rows.Item(40).RowHeight = 28.8

# --- Corrected GPS coordinates for rows 14 and 15 ---
$ws.Range("J14").Value = 40.769410000000001
$ws.Range("K14").Value = -73.996311000000006
$ws.Range("J15").Value = 40.689183999999997
$ws.Range("K15").Value = -74.044769000000002

# --- Scroll position: user had scrolled further down and selected L44 ---
[void]$ws.Range("L44").Select()
